$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A38").Value = "Nguyễn Minh Thảo"
$ws2.Range("B38").Value = "Danh mục các món ăn nên gom nhóm lại theo dạng giống menu: món chính, món phụ, món uống…"
$ws2.Range("B39").Value = "Thiếu thông tin thời gian phục vụ của nhà hàng"
$ws2.Range("A41").Value = "Nguyễn Chí Hiếu"
$ws2.Range("B41").Value = "Trang web thiết kế khó sử dụng, cần hoàn thiện các chức năng hơn, trang web cũng không thấy có gì đặc sắc"

$ws2.Range("B41").Select()
